$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.618.12'
$ws.Range('E2').Value = '  -2.45%  '
$ws.Range('D3').Value = '1.657.87'
$ws.Range('E3').Value = '  -4.23%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.03'
$ws.Range('E5').Value = '  -1.99%  '
$ws.Range('E6').Value = '  -2.30%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.15'
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  -2.23%  '
$ws.Range('E10').Value = '  -2.64%  '
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('D12').Value = '1.892.68'
$ws.Range('E12').Value = '  -4.26%  '
$ws.Range('D13').Value = '1.649.94'
$ws.Range('E13').Value = '  -4.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.15'
$ws.Range('E14').Value = '  -2.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.567'
$ws.Range('E15').Value = '  +0.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.89'
$ws.Range('E16').Value = '  -2.88%  '
$ws.Range('D17').Value = '27.605.05'
$ws.Range('E17').Value = '  -2.51%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '241.29'
$ws.Range('E18').Value = '  -2.76%  '
$ws.Range('D19').Value = '0.0₃0729'
$ws.Range('E19').Value = '  -3.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.56'
$ws.Range('E20').Value = '  -4.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('E23').Value = '  -3.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.06'
$ws.Range('E24').Value = '  -2.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.97'
$ws.Range('E25').Value = '  -2.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.22'
$ws.Range('E26').Value = '  -4.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.31'
$ws.Range('E27').Value = '  -2.47%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.111'
$ws.Range('E29').Value = '  -2.16%  '
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0503'
$ws.Range('E31').Value = '  -2.68%  '
$ws.Range('D33').Value = '1.456.21'
$ws.Range('E33').Value = '  -2.08%  '
$ws.Range('E34').Value = '  -4.94%  '
$ws.Range('E35').Value = '  -4.82%  '
$ws.Range('E36').Value = '  -1.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.924'
$ws.Range('E37').Value = '  -5.81%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.574'
$ws.Range('E38').Value = '  -4.71%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0172'
$ws.Range('E39').Value = '  -2.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '69.74'
$ws.Range('E40').Value = '  -0.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.03'
$ws.Range('E41').Value = '  -3.80%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('E43').Value = '  -4.14%  '
$ws.Range('E44').Value = '  -3.21%  '
$ws.Range('E45').Value = '  -0.32%  '
$ws.Range('D46').Value = '1.801.21'
$ws.Range('E46').Value = '  -4.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.71'
$ws.Range('E47').Value = '  -1.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.61'
$ws.Range('E48').Value = '  -2.26%  '
$ws.Range('E49').Value = '  -6.21%  '
$ws.Range('E50').Value = '  -1.27%  '
$ws.Range('E51').Value = '  -4.53%  '
